$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the 4th review row's extra columns (A4:E4), keeping only F4's
#    review text ("Exclusive info and great explanations!! bitcoin is hottt").
#    Clear() (not ClearContents) drops the cells outright so the row matches
#    the target shape of a single populated cell.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Clear()

# Row 4 had a taller custom height (13.8) left over from its old content;
# AutoFit drops back to the sheet's default row height (12.8) with no
# explicit custom-height override, same as rows 1-3.
$ws.Rows.Item(4).AutoFit()

# ---------------------------------------------------------------------------
# 2) Drop the hyperlink that used to sit on C4 (zaittomer@gmail.com), while
#    leaving the other four mail hyperlinks (C2, D2, C3, D3) intact. This
#    runtime's Hyperlinks collection only supports deleting everything at
#    once, and re-Add()-ing a hyperlink stamps the built-in "Hyperlink"
#    style onto the cell - so stash each surviving cell's original
#    formatting first and restore it (format-only paste) right after the
#    link is re-created.
# ---------------------------------------------------------------------------
$keepRanges  = @("C2", "D2", "C3", "D3")
$keepScratch = @("Z101", "Z102", "Z103", "Z104")
$keepTargets = @("mailto:budoyoni@gmail.com", "mailto:sm6502345@gmail.com", "mailto:eligitel@gmail.com", "mailto:ronenchen27@gmail.com")
$keepDisplay = @("budoyoni@gmail.com", "sm6502345@gmail.com", "eligitel@gmail.com", "ronenchen27@gmail.com")

for ($i = 0; $i -lt $keepRanges.Length; $i++) {
    $ws.Range($keepRanges[$i]).Copy($ws.Range($keepScratch[$i]))
}

$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $keepRanges.Length; $i++) {
    $ws.Hyperlinks.Add($ws.Range($keepRanges[$i]), $keepTargets[$i], "", "", $keepDisplay[$i])
}

for ($i = 0; $i -lt $keepRanges.Length; $i++) {
    $ws.Range($keepScratch[$i]).Copy()
    $ws.Range($keepRanges[$i]).PasteSpecial(-4122)
    $ws.Range($keepScratch[$i]).Clear()
}

# Adding hyperlinks registers a built-in "Hyperlink" cell style even though
# we immediately overwrote the cell formatting above; drop the now-unused
# style entry so the style table matches its pre-edit shape.
$wb.Styles.Item($wb.Styles.Count()).Delete()
